$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, preserving its original style
# (needed for numeric-looking strings in column D, which Excel would
# otherwise auto-convert to a Number cell).
function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.Value = "'" + $text
    $rng.Style = $origStyle
}

Set-TextValue "D2" "235.62"
Set-TextValue "D3" "21.76"
$ws.Range("B4").Value = "LEO"
$ws.Range("C4").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D4" "3.935"
$ws.Range("E4").Value = "3LEOLEOBestin24h"
$ws.Range("B5").Value = "HuobiToken"
$ws.Range("C5").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D5" "5.374"
$ws.Range("E5").Value = "4HuobiTokenHT"
$ws.Range("B6").Value = "Cronos"
$ws.Range("C6").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D6" "0.05578"
$ws.Range("E6").Value = "5CronosCRO"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "3.365"
$ws.Range("E7").Value = "6GateTokenGT"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "6.461"
$ws.Range("E8").Value = "7KuCoinTokenKCS"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D9" "0.8040"
$ws.Range("E9").Value = "8MXTokenMX"
$ws.Range("B10").Value = "FTXToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue "D10" "1.036"
$ws.Range("E10").Value = "9FTXTokenFTT"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue "D11" "0.1399"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.07302"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue "D13" "0.03153"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "ProBitToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D14" "0.1248"
$ws.Range("E14").Value = "13ProBitTokenPROB"
$ws.Range("B15").Value = "BitrueCoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D15" "0.02926"
$ws.Range("E15").Value = "14BitrueCoinBTR"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue "D16" "0.09243"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue "D17" "0.001659"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "MCDex"
$ws.Range("C18").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue "D18" "3.258"
$ws.Range("E18").Value = "17MCDexMCB"
$ws.Range("B19").Value = "CoinExToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D19" "0.04792"
$ws.Range("E19").Value = "18CoinExTokenCET"
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue "D20" "0.0005713"
$ws.Range("E20").Value = "19OneONE"
$ws.Range("B21").Value = "TigerCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue "D21" "0.006262"
$ws.Range("E21").Value = "20TigerCashTCH"
$ws.Range("B22").Value = "HotbitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D22" "0.005082"
$ws.Range("E22").Value = "21HotbitTokenHTB"
$ws.Range("B23").Value = "BitKan"
$ws.Range("C23").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D23" "0.001049"
$ws.Range("E23").Value = "22BitKanKAN"
$ws.Range("B24").Value = "NitroEx"
$ws.Range("C24").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D24" "0.0001501"
$ws.Range("E24").Value = "23NitroExNTX"
$ws.Range("B25").Value = "UpBots"
$ws.Range("C25").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D25" "0.0004202"
$ws.Range("E25").Value = "24UpBotsUBXT"
$ws.Range("B26").Value = "BTSEToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D26" "2.200"
$ws.Range("E26").Value = "25BTSETokenBTSE"
$ws.Range("B27").Value = "BitpandaEcosystemToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D27" "0.3286"
$ws.Range("E27").Value = "26BitpandaEcosystemTokenBEST"
Set-TextValue "D40" "0.04115"
Set-TextValue "D41" "0.007018"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003503"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1035"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue "D44" "0.008944"
Set-TextValue "D45" "0.00005445"
Set-TextValue "D46" "0.00000000751"
Set-TextValue "D48" "0.03343"
Set-TextValue "D49" "0.00002102"
